# zhongshu_wangge.xlsx 20201202 check-in update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Block 1 (创业板50（159949）/ rows 3-9) ---
$ws.Range("B4").Value = 8000
$ws.Range("D9").Value = "1.097/1.103"
$ws.Range("F9").ClearContents()
$ws.Range("G9").ClearContents()

# --- Block 2 (300ETF（510300）/ rows 12-18) ---
$ws.Range("B13").Value = 1800
$ws.Range("B17").Value = "5.026/5.074"
$ws.Range("D17").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("B18").Value = "5.120/5.154"

# --- Block 3 (科创50（588000）/ rows 21-27) ---
$ws.Range("B22").Value = 6100
$ws.Range("F27").ClearContents()
$ws.Range("G27").ClearContents()

# --- View / selection update to match the saved workbook state ---
$ws.Range("C21").Select()
